$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 216.5832213333333
$ws.Range("H2").Value = 649.749664
$ws.Range("I2").Value = 0.4331411212367192
$ws.Range("J2").Value = 0.4331411212367192
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 15.35884066666667
$ws.Range("N2").Value = 46.076522
$ws.Range("O2").Value = 0.1012042817263867
$ws.Range("P2").Value = 0.1012042817263867
$ws.Range("Q2").Value = 3326.467187532067
$ws.Range("R2").Value = 29938.20468778861
$ws.Range("S2").Value = 0.04383573606092393
$ws.Range("T2").Value = 0.04383573606092393

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 216.5832213333333
$ws.Range("H3").Value = 649.749664
$ws.Range("I3").Value = 0.4331411212367192
$ws.Range("J3").Value = 0.4331411212367192
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 50.59256466666667
$ws.Range("N3").Value = 151.777694
$ws.Range("O3").Value = 0.3333704853712116
$ws.Range("P3").Value = 0.3333704853712116
$ws.Range("Q3").Value = 10957.50063102165
$ws.Range("R3").Value = 98617.50567919482
$ws.Range("S3").Value = 0.1443964658209159
$ws.Range("T3").Value = 0.1443964658209159

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 216.5832213333333
$ws.Range("H4").Value = 649.749664
$ws.Range("I4").Value = 0.4331411212367192
$ws.Range("J4").Value = 0.4331411212367192
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 60.37715666666667
$ws.Range("N4").Value = 181.13147
$ws.Range("O4").Value = 0.397844271305776
$ws.Range("P4").Value = 0.397844271305776
$ws.Range("Q4").Value = 13076.67908581401
$ws.Range("R4").Value = 117690.1117723261
$ws.Range("S4").Value = 0.1723227137509893
$ws.Range("T4").Value = 0.1723227137509893

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 216.5832213333333
$ws.Range("H5").Value = 649.749664
$ws.Range("I5").Value = 0.4331411212367192
$ws.Range("J5").Value = 0.4331411212367192
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 25.43221733333333
$ws.Range("N5").Value = 76.296652
$ws.Range("O5").Value = 0.1675809615966257
$ws.Range("P5").Value = 0.1675809615966258
$ws.Range("Q5").Value = 5508.19155570277
$ws.Range("R5").Value = 49573.72400132493
$ws.Range("S5").Value = 0.07258620560389006
$ws.Range("T5").Value = 0.07258620560389006

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 161.954974
$ws.Range("H6").Value = 485.864922
$ws.Range("I6").Value = 0.3238910133313607
$ws.Range("J6").Value = 0.3238910133313606
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 15.35884066666667
$ws.Range("N6").Value = 46.076522
$ws.Range("O6").Value = 0.1012042817263867
$ws.Range("P6").Value = 0.1012042817263867
$ws.Range("Q6").Value = 2487.440640840142
$ws.Range("R6").Value = 22386.96576756128
$ws.Range("S6").Value = 0.03277915736183188
$ws.Range("T6").Value = 0.03277915736183188

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 161.954974
$ws.Range("H7").Value = 485.864922
$ws.Range("I7").Value = 0.3238910133313607
$ws.Range("J7").Value = 0.3238910133313606
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 50.59256466666667
$ws.Range("N7").Value = 151.777694
$ws.Range("O7").Value = 0.3333704853712116
$ws.Range("P7").Value = 0.3333704853712116
$ws.Range("Q7").Value = 8193.717495183319
$ws.Range("R7").Value = 73743.45745664986
$ws.Range("S7").Value = 0.1079757043216493
$ws.Range("T7").Value = 0.1079757043216492

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 161.954974
$ws.Range("H8").Value = 485.864922
$ws.Range("I8").Value = 0.3238910133313607
$ws.Range("J8").Value = 0.3238910133313606
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 60.37715666666667
$ws.Range("N8").Value = 181.13147
$ws.Range("O8").Value = 0.397844271305776
$ws.Range("P8").Value = 0.397844271305776
$ws.Range("Q8").Value = 9778.380838143927
$ws.Range("R8").Value = 88005.42754329534
$ws.Range("S8").Value = 0.1288581841813046
$ws.Range("T8").Value = 0.1288581841813045

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 161.954974
$ws.Range("H9").Value = 485.864922
$ws.Range("I9").Value = 0.3238910133313607
$ws.Range("J9").Value = 0.3238910133313606
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 25.43221733333333
$ws.Range("N9").Value = 76.296652
$ws.Range("O9").Value = 0.1675809615966257
$ws.Range("P9").Value = 0.1675809615966258
$ws.Range("Q9").Value = 4118.874096982348
$ws.Range("R9").Value = 37069.86687284114
$ws.Range("S9").Value = 0.05427796746657495
$ws.Range("T9").Value = 0.05427796746657495

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.4608033333333333
$ws.Range("H10").Value = 1.38241
$ws.Range("I10").Value = 0.0009215527926904059
$ws.Range("J10").Value = 0.0009215527926904059
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 15.35884066666667
$ws.Range("N10").Value = 46.076522
$ws.Range("O10").Value = 0.1012042817263867
$ws.Range("P10").Value = 0.1012042817263867
$ws.Range("Q10").Value = 7.077404975335554
$ws.Range("R10").Value = 63.69664477801999
$ws.Range("S10").Value = 0.00009326508845717824
$ws.Range("T10").Value = 0.00009326508845717825

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.4608033333333333
$ws.Range("H11").Value = 1.38241
$ws.Range("I11").Value = 0.0009215527926904059
$ws.Range("J11").Value = 0.0009215527926904059
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 50.59256466666667
$ws.Range("N11").Value = 151.777694
$ws.Range("O11").Value = 0.3333704853712116
$ws.Range("P11").Value = 0.3333704853712116
$ws.Range("Q11").Value = 23.31322244028222
$ws.Range("R11").Value = 209.81900196254
$ws.Range("S11").Value = 0.0003072185017943961
$ws.Range("T11").Value = 0.0003072185017943961

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.4608033333333333
$ws.Range("H12").Value = 1.38241
$ws.Range("I12").Value = 0.0009215527926904059
$ws.Range("J12").Value = 0.0009215527926904059
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 60.37715666666667
$ws.Range("N12").Value = 181.13147
$ws.Range("O12").Value = 0.397844271305776
$ws.Range("P12").Value = 0.397844271305776
$ws.Range("Q12").Value = 27.82199504918889
$ws.Range("R12").Value = 250.3979554427
$ws.Range("S12").Value = 0.0003666344992777174
$ws.Range("T12").Value = 0.0003666344992777174

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.4608033333333333
$ws.Range("H13").Value = 1.38241
$ws.Range("I13").Value = 0.0009215527926904059
$ws.Range("J13").Value = 0.0009215527926904059
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 25.43221733333333
$ws.Range("N13").Value = 76.296652
$ws.Range("O13").Value = 0.1675809615966257
$ws.Range("P13").Value = 0.1675809615966258
$ws.Range("Q13").Value = 11.71925052125778
$ws.Range("R13").Value = 105.47325469132
$ws.Range("S13").Value = 0.0001544347031611141
$ws.Range("T13").Value = 0.0001544347031611141

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 121.0302313333333
$ws.Range("H14").Value = 363.090694
$ws.Range("I14").Value = 0.2420463126392298
$ws.Range("J14").Value = 0.2420463126392298
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 15.35884066666667
$ws.Range("N14").Value = 46.076522
$ws.Range("O14").Value = 0.1012042817263867
$ws.Range("P14").Value = 0.1012042817263867
$ws.Range("Q14").Value = 1858.884038898474
$ws.Range("R14").Value = 16729.95635008627
$ws.Range("S14").Value = 0.02449612321517368
$ws.Range("T14").Value = 0.02449612321517368

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 121.0302313333333
$ws.Range("H15").Value = 363.090694
$ws.Range("I15").Value = 0.2420463126392298
$ws.Range("J15").Value = 0.2420463126392298
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 50.59256466666667
$ws.Range("N15").Value = 151.777694
$ws.Range("O15").Value = 0.3333704853712116
$ws.Range("P15").Value = 0.3333704853712116
$ws.Range("Q15").Value = 6123.229805353293
$ws.Range("R15").Value = 55109.06824817963
$ws.Range("S15").Value = 0.08069109672685205
$ws.Range("T15").Value = 0.08069109672685204

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 121.0302313333333
$ws.Range("H16").Value = 363.090694
$ws.Range("I16").Value = 0.2420463126392298
$ws.Range("J16").Value = 0.2420463126392298
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 60.37715666666667
$ws.Range("N16").Value = 181.13147
$ws.Range("O16").Value = 0.397844271305776
$ws.Range("P16").Value = 0.397844271305776
$ws.Range("Q16").Value = 7307.461238615576
$ws.Range("R16").Value = 65767.15114754018
$ws.Range("S16").Value = 0.09629673887420441
$ws.Range("T16").Value = 0.0962967388742044

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 121.0302313333333
$ws.Range("H17").Value = 363.090694
$ws.Range("I17").Value = 0.2420463126392298
$ws.Range("J17").Value = 0.2420463126392298
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 25.43221733333333
$ws.Range("N17").Value = 76.296652
$ws.Range("O17").Value = 0.1675809615966257
$ws.Range("P17").Value = 0.1675809615966258
$ws.Range("Q17").Value = 3078.067147172943
$ws.Range("R17").Value = 27702.60432455648
$ws.Range("S17").Value = 0.04056235382299964
$ws.Range("T17").Value = 0.04056235382299964

